$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 2019 row values (70 -> 75)
$ws.Range("B53").Value = 75
$ws.Range("C53").Value = 75

# Add 2020 row
$ws.Range("A54").Value = 2020
$ws.Range("B54").Value = 20
$ws.Range("C54").Value = 20

# Add 2021 row
$ws.Range("A55").Value = 2021
$ws.Range("B55").Value = 5
$ws.Range("C55").Value = 5

# Update selection to reflect new active cell (A56) after data entry
$ws.Range("A56").Select()
